# Fix the mislabeled sequence-diagram ID under "Confirm order sequences":
# it was a duplicate of the previous section's ID (F_SQ_SelectMeal_006)
# and should read F_SQ_ConfirmMeal_006.
$d = $word.ActiveDocument

$d.Content.Find.Execute("F_SQ_SelectMeal_006", $true, $true, $false, $false, $false, $true, 1, $false, "F_SQ_ConfirmMeal_006", 2)
